$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = "29.309.84"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.873.17"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "0.7083"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "241.86"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.07794"
$ws.Range("D9").Value = "0.3107"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "0.08377"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "1.864.30"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "5.234"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "91.26"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "0.000008387"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "6.141"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").Value = "29.319.46"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "240.31"
$ws.Range("D20").Value = "2.127.69"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "7.740"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").Value = "162.66"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "18.48"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D30").Value = "4.406"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "4.339"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").Value = "1.214"
$ws.Range("E32").Value = "  -5.84%  "
$ws.Range("D33").Value = "0.05354"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").Value = "1.941"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "0.7454"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").Value = "2.684"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").Value = "1.239.83"
$ws.Range("E39").Value = "  +6.75%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").Value = "6.501"
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("D42").Value = "0.8935"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("E43").Value = "  +4.68%  "
$ws.Range("D44").Value = "72.21"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.00000000130"
$ws.Range("E45").Value = "  +12.24%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.013.33"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.5194"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.792"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.454"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "0.4333"
$ws.Range("E51").Value = "  +0.45%  "
